# Generate Report for Handback
#
# The handback for both the zh-cn and de-de localized files has landed and
# is back in sync with en-US, so the status report is refreshed:
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#     (Overview sheet's per-locale columns, and each locale sheet's Status column).
#   - The "Latest Handback DateTime" for each locale is bumped to the time the
#     handback report was regenerated.
#   - The stale "version mismatch" error on the de-de row is cleared now that
#     the handback is current.
#   - The columns whose displayed text changed are re-fit to their new content.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns (E, F) for both rows ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# --- zh-cn sheet: Status column (C), Handback DateTime (K), Error Detail (P) ---
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("K2").Value = "2016-09-09 13:32:43"
$zhcn.Range("K3").Value = "2016-09-09 13:32:43"
$zhcn.Range("P3").Value = ""

# --- de-de sheet: Status column (C), Handback DateTime (K), Error Detail (P) ---
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus
$dede.Range("K2").Value = "2016-09-09 13:33:00"
$dede.Range("K3").Value = "2016-09-09 13:33:00"
$dede.Range("P3").Value = ""

# --- Re-fit the columns whose text just changed length ---
$overview.Range("E1:F3").EntireColumn.ColumnWidth = 29.166666666666668
$zhcn.Range("C1:C3").EntireColumn.ColumnWidth = 29.166666666666668
$dede.Range("C1:C3").EntireColumn.ColumnWidth = 29.166666666666668
$zhcn.Range("P1:P3").EntireColumn.ColumnWidth = 12.833333333333334
$dede.Range("P1:P3").EntireColumn.ColumnWidth = 12.833333333333334
